# Update cryptos list data (prices, 1h volume %, and re-ranked coin order)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price column as Text so numeric-looking strings
# (e.g. "1.00", "246.01") are preserved exactly as text, matching
# the existing inline-string storage used throughout the sheet.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '36.546.80'
$ws.Range("E2").Value = '  +3.72%  '
$ws.Range("D3").Value = '2.009.28'
$ws.Range("E3").Value = '  +6.29%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '246.01'
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").Value = '0.656'
$ws.Range("E6").Value = '  -4.44%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '44.88'
$ws.Range("E8").Value = '  +4.81%  '
$ws.Range("D9").Value = '0.362'
$ws.Range("E9").Value = '  +1.42%  '
$ws.Range("D10").Value = '56.16'
$ws.Range("E10").Value = '  +2.15%  '
$ws.Range("D11").Value = '0.0717'
$ws.Range("E11").Value = '  -3.36%  '
$ws.Range("D12").Value = '0.0990'
$ws.Range("E12").Value = '  +1.22%  '
$ws.Range("D13").Value = '14.40'
$ws.Range("E13").Value = '  +3.16%  '
$ws.Range("D14").Value = '2.294.70'
$ws.Range("E14").Value = '  +5.99%  '
$ws.Range("D15").Value = '0.801'
$ws.Range("E15").Value = '  +2.09%  '
$ws.Range("D16").Value = '2.014.40'
$ws.Range("E16").Value = '  +6.38%  '
$ws.Range("D17").Value = '4.88'
$ws.Range("E17").Value = '  -1.43%  '
$ws.Range("D18").Value = '36.681.94'
$ws.Range("E18").Value = '  +4.06%  '
$ws.Range("D19").Value = '71.08'
$ws.Range("E19").Value = '  -3.24%  '
$ws.Range("D20").Value = '0.0₃0812'
$ws.Range("E20").Value = '  -1.26%  '
$ws.Range("D21").Value = '12.95'
$ws.Range("E21").Value = '  +1.17%  '
$ws.Range("D22").Value = '234.33'
$ws.Range("E22").Value = '  -3.77%  '
$ws.Range("D23").Value = '4.98'
$ws.Range("E23").Value = '  -4.08%  '
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").Value = '2.44'
$ws.Range("E25").Value = '  -7.77%  '
$ws.Range("D26").Value = '161.90'
$ws.Range("E26").Value = '  -3.34%  '
$ws.Range("D27").Value = '1.98'
$ws.Range("E27").Value = '  -8.47%  '
$ws.Range("D28").Value = '19.62'
$ws.Range("E28").Value = '  +7.34%  '
$ws.Range("D29").Value = '8.50'
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("D30").Value = '0.122'
$ws.Range("E30").Value = '  -4.07%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '4.35'
$ws.Range("E31").Value = '  +0.61%  '
$ws.Range("B32").Value = 'Gas'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D32").Value = '20.57'
$ws.Range("E32").Value = '  +54.31%  '
$ws.Range("D33").Value = '0.0579'
$ws.Range("E33").Value = '  -2.90%  '
$ws.Range("B34").Value = 'BinanceUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '1.85'
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").Value = '4.04'
$ws.Range("E36").Value = '  -3.52%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '0.0800'
$ws.Range("E37").Value = '  +12.26%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").Value = '2.12'
$ws.Range("E38").Value = '  +8.97%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '0.844'
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '1.33'
$ws.Range("E40").Value = '  -8.76%  '
$ws.Range("D41").Value = '0.0216'
$ws.Range("E41").Value = '  -2.54%  '
$ws.Range("D42").Value = '96.24'
$ws.Range("E42").Value = '  -2.40%  '
$ws.Range("D43").Value = '2.76'
$ws.Range("E43").Value = '  +14.36%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Value = '16.09'
$ws.Range("E44").Value = '  -5.96%  '
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").Value = '1.07'
$ws.Range("E45").Value = '  -0.49%  '
$ws.Range("D46").Value = '1.299.14'
$ws.Range("E46").Value = '  -2.56%  '
$ws.Range("D47").Value = '0.0816'
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("B48").Value = 'MXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D48").Value = '2.75'
$ws.Range("E48").Value = '  +0.83%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '2.21'
$ws.Range("E49").Value = '  -6.11%  '
$ws.Range("D50").Value = '2.203.47'
$ws.Range("E50").Value = '  +6.67%  '
$ws.Range("D51").Value = '3.78'
$ws.Range("E51").Value = '  +14.17%  '
